$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macro_taxonomy")

# Make this sheet the active one (matches activeTab / tabSelected shifting
# from "Costs" to "Macro_taxonomy" in the saved file).
$ws.Activate()

# Insert a new row above row 16 (pushes the old rows 16-24 down to 17-25,
# and carries the A/B/C/D column styles of row 15 down into the new blank
# row 16, matching the target formatting exactly).
$ws.Rows("16:16").Insert()

# Row 15 (Other / Urban / MATO) proportion was split in two, 1 -> 0.5.
$ws.Range("D15").Value = 0.5

# New row 16: Other / Urban / ME+MEO/LWAL / 0.5 (the other half of the
# split MATO proportion above).
$ws.Range("A16").Value = "Other"
$ws.Range("B16").Value = "Urban"
$ws.Range("C16").Value = "ME+MEO/LWAL"
$ws.Range("D16").Value = 0.5

# Old row 24 (now shifted to row 25: Other / Rural / MATO) proportion was
# likewise split, 1 -> 0.5.
$ws.Range("D25").Value = 0.5

# Two new rows appended after it with the remaining split proportions.
$ws.Range("A26").Value = "Other"
$ws.Range("B26").Value = "Rural"
$ws.Range("C26").Value = "EWV/LN"
$ws.Range("D26").Value = 0.25

$ws.Range("A27").Value = "Other"
$ws.Range("B27").Value = "Rural"
$ws.Range("C27").Value = "ME+MEO/LWAL"
$ws.Range("D27").Value = 0.25

# Leave the cursor on D17, matching the saved selection.
$null = $ws.Range("D17").Select()
